# The deck's second slide holds a single top-level group ("Group 4") that in
# turn contains a nested group ("Group 8") wrapping a picture, two freeform
# shapes and a textbox, plus two more freeform shapes that sit directly in
# the outer group. This edit ungroups the inner "Group 8" so all six shapes
# become direct, flattened children of the (re-created) outer group.
#
# PowerPoint's object model doesn't expose a nested group directly through
# Shapes/GroupItems once it is embedded several levels deep in the same
# collection item, so we do this the same way a user would in the UI:
#   1. Ungroup the outer group -> its former members land on the slide as
#      top-level shapes (the inner group included).
#   2. Ungroup that now-top-level inner group -> its four members also
#      land on the slide as top-level shapes.
#   3. Re-select all six shapes and Group() them back together, which
#      recreates a single outer group (PowerPoint assigns it a fresh
#      id/name, same as it would interactively).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$outer = $s.Shapes.Item("Group 4")

$afterFirstUngroup = $outer.Ungroup()

$innerGroupName = $null
for ($i = 1; $i -le $afterFirstUngroup.Count; $i++) {
    $candidate = $afterFirstUngroup.Item($i)
    if ($candidate.Type -eq 6) {
        $innerGroupName = $candidate.Name
    }
}

$innerGroup = $s.Shapes.Item($innerGroupName)
$afterSecondUngroup = $innerGroup.Ungroup()

$memberNames = @()
for ($i = 1; $i -le $afterSecondUngroup.Count; $i++) {
    $memberNames += $afterSecondUngroup.Item($i).Name
}
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $name = $s.Shapes.Item($i).Name
    if (-not ($memberNames -contains $name)) {
        $memberNames += $name
    }
}

$range = $s.Shapes.Range($memberNames)
$newGroup = $range.Group()
